# Correct the misspelled unit string "microequivlanetsPerLiter" to the
# correctly spelled "microequivalentsPerLiter" wherever it appears in the
# "unit" column (column C) of Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Columns.Item(3).Replace("microequivlanetsPerLiter", "microequivalentsPerLiter", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
